$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H94").Value = 5551.25
$ws.Range("I94").Value = 4068.3333
$ws.Range("K94").Value = 4068.3333
$ws.Range("M94").Value = -3617.3333
$ws.Range("H106").Value = 6211.6875
$ws.Range("I106").Value = 4224.706
$ws.Range("J106").Value = 8463.6
$ws.Range("K106").Value = 4224.706
$ws.Range("L106").Value = 8463.6
$ws.Range("M106").Value = -3593.706
$ws.Range("N106").Value = -9725.6
$ws.Range("H125").Value = 1000.6667
$ws.Range("I125").Value = 1000.6667
$ws.Range("K125").Value = 9006.0003
$ws.Range("M125").Value = -6546.0003
$ws.Range("H137").Value = 1794387.5
$ws.Range("I137").Value = 1939.1666
$ws.Range("J137").Value = 7939924.5
$ws.Range("K137").Value = 5817.4998
$ws.Range("L137").Value = 23819773.5
$ws.Range("M137").Value = -3267.4998
$ws.Range("N137").Value = -23824873.5
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H52").Value = 17999.666
$ws.Range("J52").Value = 17999.666
$ws.Range("L52").Value = 17999.666
$ws.Range("N52").Value = -18635.666
$ws.Range("H61").Value = 4055.6875
$ws.Range("I61").Value = 2838.5833
$ws.Range("J61").Value = 7707
$ws.Range("K61").Value = 2838.5833
$ws.Range("L61").Value = 7707
$ws.Range("M61").Value = -2626.5833
$ws.Range("N61").Value = -8131
$ws.Range("H122").Value = 1635.525
$ws.Range("I122").Value = 1291.7693
$ws.Range("J122").Value = 2273.9285
$ws.Range("K122").Value = 3875.3079
$ws.Range("L122").Value = 6821.7855
$ws.Range("M122").Value = -1425.3079
$ws.Range("N122").Value = -11721.7855
$ws.Range("H136").Value = 4055.6875
$ws.Range("I136").Value = 2838.5833
$ws.Range("J136").Value = 7707
$ws.Range("K136").Value = 8515.749899999999
$ws.Range("L136").Value = 23121
$ws.Range("M136").Value = -5965.749899999999
$ws.Range("N136").Value = -28221
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H14").Value = 1680
$ws.Range("I14").Value = 1680
$ws.Range("J14").Value = 0
$ws.Range("K14").Value = 1680
$ws.Range("L14").Value = 0
$ws.Range("M14").Value = -1508
$ws.Range("N14").ClearContents()
$ws.Range("H117").Value = 30888.555
$ws.Range("I117").Value = 0
$ws.Range("J117").Value = 30888.555
$ws.Range("K117").Value = 0
$ws.Range("L117").Value = 30888.555
$ws.Range("M117").ClearContents()
$ws.Range("N117").Value = -40066.555
$ws.Range("H134").Value = 7294.2856
$ws.Range("I134").Value = 8475.875
$ws.Range("J134").Value = 5718.8335
$ws.Range("K134").Value = 25427.625
$ws.Range("L134").Value = 17156.5005
$ws.Range("M134").Value = -22892.625
$ws.Range("N134").Value = -22226.5005
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1977
$ws.Range("I31").Value = 1379.5333
$ws.Range("K31").Value = 1379.5333
$ws.Range("M31").Value = -1084.5333
$ws.Range("H34").Value = 1977
$ws.Range("I34").Value = 1379.5333
$ws.Range("K34").Value = 1379.5333
$ws.Range("M34").Value = -1177.5333
$ws.Range("H99").Value = 2220
$ws.Range("I99").Value = 2500
$ws.Range("J99").Value = 2100
$ws.Range("K99").Value = 2500
$ws.Range("L99").Value = 2100
$ws.Range("M99").Value = -1002
$ws.Range("N99").Value = -5096
$ws.Range("H105").Value = 1310.1666
$ws.Range("I105").Value = 1072.2
$ws.Range("J105").Value = 2500
$ws.Range("K105").Value = 1072.2
$ws.Range("L105").Value = 2500
$ws.Range("M105").Value = 674.8
$ws.Range("N105").Value = -5994
$ws.Range("H126").Value = 2220
$ws.Range("I126").Value = 2500
$ws.Range("J126").Value = 2100
$ws.Range("K126").Value = 7500
$ws.Range("L126").Value = 6300
$ws.Range("M126").Value = -5030
$ws.Range("N126").Value = -11240
$ws.Range("H134").Value = 2280.2222
$ws.Range("I134").Value = 2067.8
$ws.Range("J134").Value = 2887.1428
$ws.Range("K134").Value = 6203.400000000001
$ws.Range("L134").Value = 8661.428400000001
$ws.Range("M134").Value = -3668.400000000001
$ws.Range("N134").Value = -13731.4284
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 498.6857
$ws.Range("I5").Value = 334.10715
$ws.Range("J5").Value = 1157
$ws.Range("K5").Value = 1002.32145
$ws.Range("L5").Value = 3471
$ws.Range("M5").Value = -890.3214499999999
$ws.Range("N5").Value = -3695
$ws.Range("H51").Value = 3409.524
$ws.Range("J51").Value = 3409.524
$ws.Range("L51").Value = 10228.572
$ws.Range("N51").Value = -11148.572
$ws.Range("H122").Value = 1471.6
$ws.Range("I122").Value = 1063.3334
$ws.Range("J122").Value = 1573.6666
$ws.Range("K122").Value = 9570.000599999999
$ws.Range("L122").Value = 14162.9994
$ws.Range("M122").Value = -7120.000599999999
$ws.Range("N122").Value = -19062.9994
$ws.Range("H132").Value = 2344.0386
$ws.Range("I132").Value = 2097
$ws.Range("J132").Value = 3167.5
$ws.Range("K132").Value = 18873
$ws.Range("L132").Value = 28507.5
$ws.Range("M132").Value = -16343
$ws.Range("N132").Value = -33567.5
$ws.Range("H135").Value = 498.6857
$ws.Range("I135").Value = 334.10715
$ws.Range("J135").Value = 1157
$ws.Range("K135").Value = 3006.96435
$ws.Range("L135").Value = 10413
$ws.Range("M135").Value = -471.9643499999997
$ws.Range("N135").Value = -15483
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 2808.913
$ws.Range("I80").Value = 2722.5
$ws.Range("J80").Value = 3120
$ws.Range("K80").Value = 2722.5
$ws.Range("L80").Value = 3120
$ws.Range("M80").Value = -1724.5
$ws.Range("N80").Value = -5116
$ws.Range("H83").Value = 2808.913
$ws.Range("I83").Value = 2722.5
$ws.Range("J83").Value = 3120
$ws.Range("K83").Value = 13612.5
$ws.Range("L83").Value = 15600
$ws.Range("M83").Value = -8620.5
$ws.Range("N83").Value = -25584
$ws.Range("H122").Value = 60384
$ws.Range("I122").Value = 92355.17999999999
$ws.Range("J122").Value = 1770.1666
$ws.Range("K122").Value = 277065.54
$ws.Range("L122").Value = 5310.4998
$ws.Range("M122").Value = -274615.54
$ws.Range("N122").Value = -10210.4998
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2737.3333
$ws.Range("J7").Value = 2852.8
$ws.Range("L7").Value = 2852.8
$ws.Range("N7").Value = -3076.8
$ws.Range("H22").Value = 993.3333
$ws.Range("J22").Value = 992
$ws.Range("L22").Value = 992
$ws.Range("N22").Value = -1582
$ws.Range("H27").Value = 993.3333
$ws.Range("J27").Value = 992
$ws.Range("L27").Value = 992
$ws.Range("N27").Value = -1206
$ws.Range("H68").Value = 15812.5
$ws.Range("I68").Value = 52250
$ws.Range("J68").Value = 3666.6667
$ws.Range("K68").Value = 52250
$ws.Range("L68").Value = 3666.6667
$ws.Range("M68").Value = -51501
$ws.Range("N68").Value = -5164.6667
$ws.Range("H71").Value = 15812.5
$ws.Range("I71").Value = 52250
$ws.Range("J71").Value = 3666.6667
$ws.Range("K71").Value = 261250
$ws.Range("L71").Value = 18333.3335
$ws.Range("M71").Value = -257506
$ws.Range("N71").Value = -25821.3335
$ws.Range("H115").Value = 0
$ws.Range("J115").Value = 0
$ws.Range("L115").Value = 0
$ws.Range("N115").ClearContents()
$ws.Range("H122").Value = 2809.9688
$ws.Range("I122").Value = 2256.5293
$ws.Range("K122").Value = 6769.5879
$ws.Range("M122").Value = -4319.5879
$ws.Range("H126").Value = 2737.3333
$ws.Range("J126").Value = 2852.8
$ws.Range("L126").Value = 8558.400000000001
$ws.Range("N126").Value = -13498.4
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 1213.8928
$ws.Range("I126").Value = 683.125
$ws.Range("K126").Value = 2049.375
$ws.Range("M126").Value = 420.625
$ws.Range("H132").Value = 3681.9033
$ws.Range("I132").Value = 3634.9048
$ws.Range("K132").Value = 10904.7144
$ws.Range("M132").Value = -8374.714399999999
